# Auto-generated edit script: updates market-board derived profit
# calculations (currentAveragePrice*, LevePrice*, LeveProfit* columns)
# across all 8 job sheets, reflecting refreshed Universalis price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3006193.8
$ws.Range("I15").Value = 3006193.8
$ws.Range("K15").Value = 9018581.399999999
$ws.Range("M15").Value = -9018412.399999999
$ws.Range("H40").Value = 4051.1667
$ws.Range("I40").Value = 3365.7144
$ws.Range("K40").Value = 3365.7144
$ws.Range("M40").Value = -3190.7144
$ws.Range("H75").Value = 29999
$ws.Range("J75").Value = 29999
$ws.Range("L75").Value = 29999
$ws.Range("N75").Value = -31871
$ws.Range("H76").Value = 5885.6924
$ws.Range("I76").Value = 4502
$ws.Range("J76").Value = 7500
$ws.Range("K76").Value = 4502
$ws.Range("L76").Value = 7500
$ws.Range("M76").Value = -4187
$ws.Range("N76").Value = -8130
$ws.Range("H78").Value = 29999
$ws.Range("J78").Value = 29999
$ws.Range("L78").Value = 89997
$ws.Range("N78").Value = -99357
$ws.Range("H79").Value = 5885.6924
$ws.Range("I79").Value = 4502
$ws.Range("J79").Value = 7500
$ws.Range("K79").Value = 4502
$ws.Range("L79").Value = 7500
$ws.Range("M79").Value = -3410
$ws.Range("N79").Value = -9684
$ws.Range("H113").Value = 35716650
$ws.Range("I113").Value = 12502136
$ws.Range("K113").Value = 12502136
$ws.Range("M113").Value = -12498882
$ws.Range("H138").Value = 1733.17
$ws.Range("I138").Value = 1681
$ws.Range("J138").Value = 1738.9667
$ws.Range("K138").Value = 5043
$ws.Range("L138").Value = 5216.9001
$ws.Range("M138").Value = 97
$ws.Range("N138").Value = -15496.9001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 603.28
$ws.Range("I2").Value = 489.61905
$ws.Range("K2").Value = 489.61905
$ws.Range("M2").Value = -376.61905
$ws.Range("H7").Value = 56997.5
$ws.Range("I7").Value = 35000
$ws.Range("J7").Value = 78995
$ws.Range("K7").Value = 35000
$ws.Range("L7").Value = 78995
$ws.Range("M7").Value = -34886
$ws.Range("N7").Value = -79223
$ws.Range("H19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("N19").Value = -1458
$ws.Range("H63").Value = 6497.9165
$ws.Range("I63").Value = 2392.6
$ws.Range("K63").Value = 2392.6
$ws.Range("M63").Value = -1706.6
$ws.Range("H66").Value = 6497.9165
$ws.Range("I66").Value = 2392.6
$ws.Range("K66").Value = 11963
$ws.Range("M66").Value = -8531
$ws.Range("H116").Value = 603.28
$ws.Range("I116").Value = 489.61905
$ws.Range("K116").Value = 489.61905
$ws.Range("M116").Value = 1804.38095
$ws.Range("H132").Value = 13691.5
$ws.Range("I132").Value = 5161
$ws.Range("K132").Value = 15483
$ws.Range("M132").Value = -12953

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 603.28
$ws.Range("I3").Value = 489.61905
$ws.Range("K3").Value = 489.61905
$ws.Range("M3").Value = -375.61905
$ws.Range("H14").Value = 554
$ws.Range("I14").Value = 554
$ws.Range("K14").Value = 554
$ws.Range("M14").Value = -382
$ws.Range("H20").Value = 6935.8
$ws.Range("I20").Value = 6920.375
$ws.Range("K20").Value = 6920.375
$ws.Range("M20").Value = -6673.375
$ws.Range("H26").Value = 19471
$ws.Range("I26").Value = 19471
$ws.Range("K26").Value = 19471
$ws.Range("M26").Value = -19179
$ws.Range("H40").Value = 44495
$ws.Range("J40").Value = 44495
$ws.Range("L40").Value = 44495
$ws.Range("N40").Value = -45025

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 951.8570999999999
$ws.Range("I22").Value = 621.8182
$ws.Range("J22").Value = 2162
$ws.Range("K22").Value = 621.8182
$ws.Range("L22").Value = 2162
$ws.Range("M22").Value = -271.8182
$ws.Range("N22").Value = -2862
$ws.Range("H132").Value = 1854.6207
$ws.Range("I132").Value = 1829.0741
$ws.Range("K132").Value = 5487.2223
$ws.Range("M132").Value = -2957.2223
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 36462380
$ws.Range("I4").Value = 24400216
$ws.Range("K4").Value = 73200648
$ws.Range("M4").Value = -73200536
$ws.Range("H5").Value = 2052.0715
$ws.Range("I5").Value = 2052.0715
$ws.Range("K5").Value = 6156.2145
$ws.Range("M5").Value = -6044.2145
$ws.Range("H7").Value = 1217.875
$ws.Range("I7").Value = 360.5
$ws.Range("K7").Value = 1081.5
$ws.Range("M7").Value = -969.5
$ws.Range("H20").Value = 8900
$ws.Range("J20").Value = 8900
$ws.Range("L20").Value = 26700
$ws.Range("N20").Value = -27154
$ws.Range("H92").Value = 308.66666
$ws.Range("I92").Value = 299
$ws.Range("K92").Value = 897
$ws.Range("M92").Value = 351
$ws.Range("H118").Value = 6287.25
$ws.Range("I118").Value = 2299.3333
$ws.Range("K118").Value = 6897.999899999999
$ws.Range("M118").Value = -5654.999899999999
$ws.Range("H135").Value = 2052.0715
$ws.Range("I135").Value = 2052.0715
$ws.Range("K135").Value = 18468.6435
$ws.Range("M135").Value = -15933.6435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3696.9412
$ws.Range("I126").Value = 3128.3333
$ws.Range("J126").Value = 4336.625
$ws.Range("K126").Value = 9384.999899999999
$ws.Range("L126").Value = 13009.875
$ws.Range("M126").Value = -6914.999899999999
$ws.Range("N126").Value = -17949.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 97617.82000000001
$ws.Range("I7").Value = 3755.4
$ws.Range("K7").Value = 3755.4
$ws.Range("M7").Value = -3643.4
$ws.Range("H22").Value = 1583.6666
$ws.Range("H27").Value = 1583.6666
$ws.Range("H46").Value = 2929.0588
$ws.Range("J46").Value = 4316.4
$ws.Range("L46").Value = 4316.4
$ws.Range("N46").Value = -4692.4
$ws.Range("H61").Value = 1933
$ws.Range("I61").Value = 2178
$ws.Range("K61").Value = 2178
$ws.Range("M61").Value = -1976
$ws.Range("H93").Value = 90910670
$ws.Range("I93").Value = 250000830
$ws.Range("K93").Value = 250000830
$ws.Range("M93").Value = -249999582
$ws.Range("H113").Value = 1933
$ws.Range("I113").Value = 2178
$ws.Range("K113").Value = 2178
$ws.Range("M113").Value = -8
$ws.Range("H121").Value = 99995
$ws.Range("J121").Value = 99995
$ws.Range("L121").Value = 99995
$ws.Range("N121").Value = -103489
$ws.Range("H123").Value = 54656.668
$ws.Range("J123").Value = 54656.668
$ws.Range("L123").Value = 54656.668
$ws.Range("N123").Value = -64456.668
$ws.Range("H126").Value = 97617.82000000001
$ws.Range("I126").Value = 3755.4
$ws.Range("K126").Value = 11266.2
$ws.Range("M126").Value = -8796.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 11905402
$ws.Range("I107").Value = 17242094
$ws.Range("J107").Value = 474.07693
$ws.Range("K107").Value = 51726282
$ws.Range("L107").Value = 1422.23079
$ws.Range("M107").Value = -51724362
$ws.Range("N107").Value = -5262.23079
$ws.Range("H132").Value = 419107.75
$ws.Range("I132").Value = 2208.2273
$ws.Range("K132").Value = 6624.6819
$ws.Range("M132").Value = -4094.6819
$ws.Range("H141").Value = 67000
$ws.Range("J141").Value = 67000
$ws.Range("L141").Value = 67000
$ws.Range("N141").Value = -77360

